$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D4"  = -7.653400000000001
    "A8"  = -22.30630000000001
    "A10" = -22.10229999999999
    "A12" = -21.4483
    "D12" = -5.961499999999997
    "D15" = -8.247399999999995
    "D17" = -8.364099999999992
    "A18" = -22.08370000000001
    "D26" = -7.331000000000002
    "D27" = -7.8538
    "D28" = -7.788099999999998
    "A37" = -20.02139999999998
    "D37" = -7.852799999999997
    "D47" = -7.715599999999998
    "A55" = -21.94080000000001
    "D65" = -7.877300000000004
    "A68" = -21.36739999999999
    "D73" = -8.197199999999993
    "A77" = -19.98309999999999
    "A78" = -19.71289999999998
    "A81" = -22.0777
    "A82" = -21.73870000000001
    "D84" = -7.6134
    "D85" = -8.589299999999998
    "D93" = -6.749599999999991
    "D95" = -7.600900000000004
    "D98" = -7.487600000000001
    "D99" = -7.929300000000003
    "D101" = -8.597599999999991
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
